# Apply "custom accuracy" rounding to row 5 (2 decimal places) and drop the
# extra data row (row 6), shrinking the sheet's used range from A1:AH6 to A1:AH5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Round row 5 measurement values (columns B:AH) to 2 decimal places ---
$ws.Range("B5").Value = 5.7
$ws.Range("C5").Value = 4.06
$ws.Range("D5").Value = 0.74
$ws.Range("E5").Value = 12.58
$ws.Range("F5").Value = 9.75
$ws.Range("G5").Value = 4.42
$ws.Range("H5").Value = 19.94
$ws.Range("I5").Value = 6.98
$ws.Range("J5").Value = 2.99
$ws.Range("K5").Value = 4.26
$ws.Range("L5").Value = 5.01
$ws.Range("M5").Value = 5.37
$ws.Range("N5").Value = 1.45
$ws.Range("O5").Value = 4.51
$ws.Range("P5").Value = 6.33
$ws.Range("Q5").Value = 4.03
$ws.Range("R5").Value = 0.69
$ws.Range("S5").Value = 0.4
$ws.Range("T5").Value = 61.36
$ws.Range("U5").Value = 12.73
$ws.Range("V5").Value = 4.16
$ws.Range("W5").Value = 8.33
$ws.Range("X5").Value = 4.32
$ws.Range("Y5").Value = 0.91
$ws.Range("Z5").Value = 9.39
$ws.Range("AA5").Value = 3.68
$ws.Range("AB5").Value = 3.37
$ws.Range("AC5").Value = 3.96
$ws.Range("AD5").Value = 5.15
$ws.Range("AE5").Value = 0.54
$ws.Range("AF5").Value = 18.4
$ws.Range("AG5").Value = 2.25
$ws.Range("AH5").Value = 5.21

# --- Remove the now-superfluous row 6 entirely (data + row itself) ---
$ws.Rows(6).Delete()
